# Final_Db_Design.xlsx — DB_Design sheet
# Commit: "+ Added search box control; + Added order screen;"
#
# Net effect on the "Ca làm việc" (work-shift) mini table in column I:
#   I16 becomes a new "Mã nhân viên" entry, bold+italic, boxed with a thin
#   border on left/top/right (no bottom), and the existing "Ca số 1/2/3"
#   rows shift down by one (I17/I18/I19); I15's "Ngày" label loses its bold
#   weight; the selection moves to I15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the active selection to I15 (was J14)
$ws.Range("I15").Select() | Out-Null

# I15: "Ngày" — de-emphasize from bold to plain
$i15 = $ws.Range("I15")
$i15.Value2 = "Ngày"
$i15.Font.Bold = $false
$i15.Font.Italic = $false

# I16: new "Mã nhân viên" entry — bold + italic, boxed border (no bottom)
$i16 = $ws.Range("I16")
$i16.Value2 = "Mã nhân viên"
$i16.Font.Bold = $true
$i16.Font.Italic = $true
$i16.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$i16.Borders.Item(7).Weight = 2      # xlThin
$i16.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$i16.Borders.Item(8).Weight = 2      # xlThin
$i16.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$i16.Borders.Item(10).Weight = 2     # xlThin
$i16.Borders.Item(9).LineStyle = -4142  # xlEdgeBottom -> xlLineStyleNone

# I17: shifted value "Ca số 2" -> "Ca số 1", stays plain/no border
$i17 = $ws.Range("I17")
$i17.Value2 = "Ca số 1"
$i17.Font.Bold = $false
$i17.Font.Italic = $false

# I18: shifted value "Ca số 3" -> "Ca số 2", stays plain/no border
$i18 = $ws.Range("I18")
$i18.Value2 = "Ca số 2"
$i18.Font.Bold = $false
$i18.Font.Italic = $false

# I19: previously "Mã nhân viên" (bold+italic) -> now "Ca số 3", plain
$i19 = $ws.Range("I19")
$i19.Value2 = "Ca số 3"
$i19.Font.Bold = $false
$i19.Font.Italic = $false

# E17 keeps the same text ("Mã chức vụ"); its style index shifts in the
# canonical file only because of the new style record, formatting unchanged.
